$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dates (stored as plain text, not real Excel dates)
$ws.Range("M2").Value = "2020-12-24 00:00:00"
$ws.Range("N2").Value = "2017-12-31 00:00:00"

# Updated numeric figures
$ws.Range("O2").Value = 393797078.78
$ws.Range("P2").Value = 121527467.58
$ws.Range("Q2").Value = 10288960.37
$ws.Range("R2").Value = ""
$ws.Range("S2").Value = 56985757.13
$ws.Range("T2").Value = ""
$ws.Range("U2").Value = 94749305.31999999
$ws.Range("V2").Value = ""
$ws.Range("W2").Value = 173026372.18
$ws.Range("X2").Value = 91015223.28
$ws.Range("Y2").Value = ""
$ws.Range("Z2").Value = 8689836.07
$ws.Range("AA2").Value = ""
$ws.Range("AB2").Value = 220770706.6
$ws.Range("AC2").Value = ""
$ws.Range("AD2").Value = ""
$ws.Range("AE2").Value = ""
$ws.Range("AF2").Value = 121.6439117345
$ws.Range("AG2").Value = 43.9379521849
